# Apply Briefing Final - Cliente.docx content updates (fills in the form answers
# and refreshes the submission timestamp) by locating each field's bold label
# paragraph and rewriting the answer text that follows the line break.
$d = $word.ActiveDocument

$updates = @(
    @{ Label = '__DATE__'; New = 'Data de envio: 19/11/2025, 14:29:28' },
    @{ Label = 'Nome completo da empresa:'; New = ' Maglietta' },
    @{ Label = 'Outros nomes/apelidos:'; New = ' Não' },
    @{ Label = 'Definição do negócio:'; New = ' Maglietta — Confecção especializada em malharia e fitness.
Atendemos pessoas, marcas e empresas que buscam uniformes, peças para eventos, viagens, aniversários ou private label, sempre com alto padrão de qualidade. Desenvolvemos produtos exclusivos, personalizados com seu logo ou estampa, utilizando materiais premium e acabamentos diferenciados — garantindo identidade, conforto e excelência em cada detalhe.' },
    @{ Label = 'CNPJ:'; New = ' 17.540.180/0001-05' },
    @{ Label = 'Categoria do negócio:'; New = ' Confecção Personalizada' },
    @{ Label = 'Endereço completo:'; New = ' R. Pedregulhos, 48 - Chácara da Barra
Campinas - SP, 13090-716' },
    @{ Label = 'Identificação na fachada:'; New = ' Não há identificação' },
    @{ Label = 'Áreas de atendimento:'; New = ' Todo o Brasil' },
    @{ Label = 'Raio de atendimento:'; New = ' Não informado' },
    @{ Label = 'Horários de funcionamento:'; New = ' Segunda a Sexta, das 8h às 18h' },
    @{ Label = 'Horários especiais:'; New = ' Nao' },
    @{ Label = 'História da criação:'; New = ' Em 2016' },
    @{ Label = 'Experiência no ramo:'; New = ' 11 anos' },
    @{ Label = 'Certificações:'; New = ' nao' },
    @{ Label = 'Prêmios e reconhecimentos:'; New = ' nao' },
    @{ Label = 'Quantidade de clientes:'; New = ' Mais de 1000' },
    @{ Label = 'Lista de produtos/serviços:'; New = ' Camiseta em diversos tecidos: 100% algodão sustentável com certificado BCI, dry 100% poliamida com proteção UV, dry 100% poliéster com proteção UV, PV Anti-pilling, Modal, Tech Modal
Polo em piquet
Moletom flanelado 2 e 3 cabos
Avental em oxford
Corta vento
Bermuda Masculina
Calça unissex
Legging
Top
Bermuda Feminina
Saia Fitness
Bolsas
Brindes
Ecobags' },
    @{ Label = 'Carro-chefe:'; New = ' Nosso carro chefe são as camisetas em malha 100% algodão e 100% poliamida.' },
    @{ Label = 'Produtos sazonais:'; New = ' Não' },
    @{ Label = 'Diferencial da concorrência:'; New = ' Trabalhamos com matéria prima e acabamentos premium para uma peça bacana e com durabilidade. Além de trabalharmos com algumas peças já em estoque produzidas, para conseguir atender mais rápido o cliente e apenas estampar. Conseguimos fazer entrega expressa.' },
    @{ Label = 'Marcas parceiras:'; New = ' Marca própria' },
    @{ Label = 'Faixa de preço:'; New = ' Médio' },
    @{ Label = 'Garantia:'; New = ' Não informado' },
    @{ Label = 'Canais de compra:'; New = ' Whastapp 19997958700' },
    @{ Label = 'Contato preferencial:'; New = ' Whastapp' },
    @{ Label = 'Formas de pagamento:'; New = ' Pix e cartão de crédito (taxas extras)' },
    @{ Label = 'Valor mínimo:'; New = ' Não' },
    @{ Label = 'Tempo de entrega:'; New = ' Conseguimos entregar super rápido o que já temos em estoque, porém se é alguma cor ou material que trabalhamos sob demanda, o prazo máximo é de até 30 dias (geralmente entregamos em 20), mas deixamos uma folga para épocas com muitos pedidos / eventuais contratempos.' },
    @{ Label = 'Instruções especiais:'; New = ' Não informado' },
    @{ Label = 'Como clientes procuram:'; New = ' Camiseta em Campinas
Camiseta personalizada
Uniforme
Camiseta com qualidade' },
    @{ Label = 'Palavras-chave desejadas:'; New = ' Camiseta
Moletom
Uniforme premium
Camiseta de algodão
Camiseta fitness
Camiseta sem pedido minimo
Uniforme para academia
Uniforme de qualidade
Camiseta para viagem' },
    @{ Label = 'Termos a evitar:'; New = ' Camiseta promocional' },
    @{ Label = 'Concorrentes diretos:'; New = ' Kan House
Espeta
Articulania
FG Confecacao
' },
    @{ Label = 'O que admira nos concorrentes:'; New = ' Com excessao da FG, a qualidade.
FG as camisetas são de qualidade inferior' },
    @{ Label = 'O que faz melhor:'; New = ' Entrega expressa e qualidade' },
    @{ Label = 'Melhor presença digital:'; New = ' Kan House' },
    @{ Label = 'Benchmark externo:'; New = ' Não informado' },
    @{ Label = 'Tipos de fotos disponíveis:'; New = ' Produtos' },
    @{ Label = 'Possui vídeos:'; New = ' Não' },
    @{ Label = 'Redes sociais:'; New = ' @maglietta_oficial
FB: magliettaoficial' },
    @{ Label = 'Site próprio:'; New = ' Nao' },
    @{ Label = 'Plataformas presentes:'; New = ' -' },
    @{ Label = 'Integrar plataformas no GMB:'; New = ' Sim' },
    @{ Label = 'Atributos do negócio:'; New = ' Estacionamento, Wi-Fi, aceita pets' },
    @{ Label = 'Características do ambiente:'; New = ' Amplo' },
    @{ Label = 'Público-alvo:'; New = ' Empresas' },
    @{ Label = 'Avaliações online recebidas:'; New = ' Não' },
    @{ Label = 'Onde recebeu avaliações:'; New = ' Não informado' },
    @{ Label = 'Estratégia para avaliações negativas:'; New = ' Não informado' },
    @{ Label = 'Estratégia para solicitar avaliações:'; New = ' Não' },
    @{ Label = 'Objetivo principal:'; New = ' Gerar mais leads' },
    @{ Label = 'Google Ads futuro:'; New = ' Talvez' },
    @{ Label = 'Responsável pela gestão:'; New = ' Rafaela 19997958700' },
    @{ Label = 'Tentativa anterior GMB:'; New = ' Não informado' },
    @{ Label = 'Informações a ocultar:'; New = ' nao' },
    @{ Label = 'Restrições legais:'; New = ' nao' },
    @{ Label = 'Produtos restritos pelo Google:'; New = ' nao' },
    @{ Label = 'Problemas anteriores com Google:'; New = ' Não' },
    @{ Label = 'Detalhes importantes:'; New = ' nao' },
    @{ Label = 'Maior expectativa:'; New = ' Melhorar a visibilidade online e aumento de vendas' },
    @{ Label = 'Orçamento para melhorias:'; New = ' Produção completa' },
)

$idx = 0
foreach ($p in $d.Paragraphs) {
    if ($idx -ge $updates.Count) { break }
    $u = $updates[$idx]
    $label = $u.Label
    $text = $p.Range.Text
    if ($label -eq '__DATE__') {
        if ($text.StartsWith('Data de envio:')) {
            $p.Range.Text = $u.New
            $idx = $idx + 1
        }
    } else {
        if ($text.StartsWith($label)) {
            $pStart = $p.Range.Start
            $pEnd = $p.Range.End
            $labelLen = $label.Length
            $valStart = $pStart + $labelLen + 1
            $rng = $d.Range($valStart, $pEnd)
            $rng.Text = $u.New
            $idx = $idx + 1
        }
    }
}

Write-Host "Applied $idx of $($updates.Count) updates."
if ($idx -ne $updates.Count) {
    throw "Not all updates were applied: $idx / $($updates.Count)"
}
